$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "hello"
Write-Output $ws.Range("A1").Value2
